# MAJ Fiches Bron pour Intégration Invictus
#
# 1) Refresh the cached "datetimeFigureOut" date field text, everywhere it
#    appears on the slide-master's custom layouts (the date placeholder on
#    each layout), from 27/02/2020 -> 06/03/2020.
# 2) On slide 1, append two new bullet paragraphs to the "Rectangle 20"
#    shape (an empty bullet line, then a new instruction line).

$p = $ppt.ActivePresentation

# --- 1) Update the cached date placeholders on every slide layout ---------
$oldDate = "27/02/2020"
$newDate = "06/03/2020"

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shape = $layout.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                    $shape.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# --- 2) Add the two new bullet paragraphs to slide 1 / "Rectangle 20" -----
$slide = $p.Slides.Item(1)
$target = $null
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $sh = $slide.Shapes.Item($k)
    if ($sh.Name -eq "Rectangle 20") {
        $target = $sh
    }
}

$newText = "`r`rD" + [char]0x00E9 + "finir un Responsable pour le Journal d" + [char]0x2019 + "Int" + [char]0x00E9 + "gration en d" + [char]0x00E9 + "but de s" + [char]0x00E9 + "ance, pour qu" + [char]0x2019 + "il le remplisse en fin de s" + [char]0x00E9 + "ance"

[void]$target.TextFrame.TextRange.InsertAfter($newText)
